$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.659958333333333
$ws.Range("H2").Value = 10.979875
$ws.Range("I2").Value = 0.4781132044744068
$ws.Range("J2").Value = 0.4781132044744067
$ws.Range("M2").Value = 3.832616
$ws.Range("N2").Value = 11.497848
$ws.Range("O2").Value = 0.264217765931355
$ws.Range("P2").Value = 0.264217765931355
$ws.Range("Q2").Value = 14.02721486766667
$ws.Range("R2").Value = 126.244933809
$ws.Range("S2").Value = 0.1263260027485089
$ws.Range("T2").Value = 0.1263260027485089

$ws.Range("G3").Value = 3.659958333333333
$ws.Range("H3").Value = 10.979875
$ws.Range("I3").Value = 0.4781132044744068
$ws.Range("J3").Value = 0.4781132044744067
$ws.Range("O3").Value = 0.6031799756961482
$ws.Range("P3").Value = 0.6031799756961482
$ws.Range("Q3").Value = 32.02258218004167
$ws.Range("R3").Value = 288.203239620375
$ws.Range("S3").Value = 0.2883883110548802
$ws.Range("T3").Value = 0.2883883110548802

$ws.Range("G4").Value = 3.659958333333333
$ws.Range("H4").Value = 10.979875
$ws.Range("I4").Value = 0.4781132044744068
$ws.Range("J4").Value = 0.4781132044744067
$ws.Range("M4").Value = 1.255882
$ws.Range("N4").Value = 3.767646
$ws.Range("O4").Value = 0.08657959375878042
$ws.Range("P4").Value = 0.08657959375878042
$ws.Range("Q4").Value = 4.596475791583333
$ws.Range("R4").Value = 41.36828212425
$ws.Range("S4").Value = 0.04139484701410286
$ws.Range("T4").Value = 0.04139484701410285

$ws.Range("G5").Value = 3.659958333333333
$ws.Range("H5").Value = 10.979875
$ws.Range("I5").Value = 0.4781132044744068
$ws.Range("J5").Value = 0.4781132044744067
$ws.Range("M5").Value = 0.6675826666666667
$ws.Range("N5").Value = 2.002748
$ws.Range("O5").Value = 0.04602266461371635
$ws.Range("P5").Value = 0.04602266461371635
$ws.Range("Q5").Value = 2.443324744055555
$ws.Range("R5").Value = 21.9899226965
$ws.Range("S5").Value = 0.02200404365691481
$ws.Range("T5").Value = 0.02200404365691481

$ws.Range("I6").Value = 0.3193330932870009
$ws.Range("J6").Value = 0.3193330932870008
$ws.Range("M6").Value = 3.832616
$ws.Range("N6").Value = 11.497848
$ws.Range("O6").Value = 0.264217765931355
$ws.Range("P6").Value = 0.264217765931355
$ws.Range("Q6").Value = 9.368814481536001
$ws.Range("R6").Value = 84.31933033382401
$ws.Range("S6").Value = 0.08437347649624036
$ws.Range("T6").Value = 0.08437347649624034

$ws.Range("I7").Value = 0.3193330932870009
$ws.Range("J7").Value = 0.3193330932870008
$ws.Range("O7").Value = 0.6031799756961482
$ws.Range("P7").Value = 0.6031799756961482
$ws.Range("S7").Value = 0.192615327447829
$ws.Range("T7").Value = 0.192615327447829

$ws.Range("I8").Value = 0.3193330932870009
$ws.Range("J8").Value = 0.3193330932870008
$ws.Range("M8").Value = 1.255882
$ws.Range("N8").Value = 3.767646
$ws.Range("O8").Value = 0.08657959375878042
$ws.Range("P8").Value = 0.08657959375878042
$ws.Range("Q8").Value = 3.069998525472
$ws.Range("R8").Value = 27.629986729248
$ws.Range("S8").Value = 0.02764772949052327
$ws.Range("T8").Value = 0.02764772949052326

$ws.Range("I9").Value = 0.3193330932870009
$ws.Range("J9").Value = 0.3193330932870008
$ws.Range("M9").Value = 0.6675826666666667
$ws.Range("N9").Value = 2.002748
$ws.Range("O9").Value = 0.04602266461371635
$ws.Range("P9").Value = 0.04602266461371635
$ws.Range("Q9").Value = 1.631903158336
$ws.Range("R9").Value = 14.687128425024
$ws.Range("S9").Value = 0.01469655985240824
$ws.Range("T9").Value = 0.01469655985240824

$ws.Range("G10").Value = 1.388093333333333
$ws.Range("H10").Value = 4.16428
$ws.Range("I10").Value = 0.1813315046964271
$ws.Range("J10").Value = 0.1813315046964271
$ws.Range("M10").Value = 3.832616
$ws.Range("N10").Value = 11.497848
$ws.Range("O10").Value = 0.264217765931355
$ws.Range("P10").Value = 0.264217765931355
$ws.Range("Q10").Value = 5.320028718826666
$ws.Range("R10").Value = 47.88025846944
$ws.Range("S10").Value = 0.04791100506386099
$ws.Range("T10").Value = 0.04791100506386098

$ws.Range("G11").Value = 1.388093333333333
$ws.Range("H11").Value = 4.16428
$ws.Range("I11").Value = 0.1813315046964271
$ws.Range("J11").Value = 0.1813315046964271
$ws.Range("O11").Value = 0.6031799756961482
$ws.Range("P11").Value = 0.6031799756961482
$ws.Range("Q11").Value = 12.14503794630667
$ws.Range("R11").Value = 109.30534151676
$ws.Range("S11").Value = 0.1093755325957369
$ws.Range("T11").Value = 0.1093755325957369

$ws.Range("G12").Value = 1.388093333333333
$ws.Range("H12").Value = 4.16428
$ws.Range("I12").Value = 0.1813315046964271
$ws.Range("J12").Value = 0.1813315046964271
$ws.Range("M12").Value = 1.255882
$ws.Range("N12").Value = 3.767646
$ws.Range("O12").Value = 0.08657959375878042
$ws.Range("P12").Value = 0.08657959375878042
$ws.Range("Q12").Value = 1.743281431653333
$ws.Range("R12").Value = 15.68953288488
$ws.Range("S12").Value = 0.01569960801228504
$ws.Range("T12").Value = 0.01569960801228504

$ws.Range("G13").Value = 1.388093333333333
$ws.Range("H13").Value = 4.16428
$ws.Range("I13").Value = 0.1813315046964271
$ws.Range("J13").Value = 0.1813315046964271
$ws.Range("M13").Value = 0.6675826666666667
$ws.Range("N13").Value = 2.002748
$ws.Range("O13").Value = 0.04602266461371635
$ws.Range("P13").Value = 0.04602266461371635
$ws.Range("Q13").Value = 0.9266670490488887
$ws.Range("R13").Value = 8.340003441439999
$ws.Range("S13").Value = 0.008345359024544197
$ws.Range("T13").Value = 0.008345359024544195

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.162456
$ws.Range("H14").Value = 0.487368
$ws.Range("I14").Value = 0.02122219754216535
$ws.Range("J14").Value = 0.02122219754216534
$ws.Range("M14").Value = 3.832616
$ws.Range("N14").Value = 11.497848
$ws.Range("O14").Value = 0.264217765931355
$ws.Range("P14").Value = 0.264217765931355
$ws.Range("Q14").Value = 0.6226314648960001
$ws.Range("R14").Value = 5.603683184064001
$ws.Range("S14").Value = 0.005607281622744821
$ws.Range("T14").Value = 0.005607281622744821

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.162456
$ws.Range("H15").Value = 0.487368
$ws.Range("I15").Value = 0.02122219754216535
$ws.Range("J15").Value = 0.02122219754216534
$ws.Range("O15").Value = 0.6031799756961482
$ws.Range("P15").Value = 0.6031799756961482
$ws.Range("Q15").Value = 1.421398862184
$ws.Range("R15").Value = 12.792589759656
$ws.Range("S15").Value = 0.01280080459770215
$ws.Range("T15").Value = 0.01280080459770215

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.162456
$ws.Range("H16").Value = 0.487368
$ws.Range("I16").Value = 0.02122219754216535
$ws.Range("J16").Value = 0.02122219754216534
$ws.Range("M16").Value = 1.255882
$ws.Range("N16").Value = 3.767646
$ws.Range("O16").Value = 0.08657959375878042
$ws.Range("P16").Value = 0.08657959375878042
$ws.Range("Q16").Value = 0.204025566192
$ws.Range("R16").Value = 1.836230095728
$ws.Range("S16").Value = 0.001837409241869264
$ws.Range("T16").Value = 0.001837409241869264

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.162456
$ws.Range("H17").Value = 0.487368
$ws.Range("I17").Value = 0.02122219754216535
$ws.Range("J17").Value = 0.02122219754216534
$ws.Range("M17").Value = 0.6675826666666667
$ws.Range("N17").Value = 2.002748
$ws.Range("O17").Value = 0.04602266461371635
$ws.Range("P17").Value = 0.04602266461371635
$ws.Range("Q17").Value = 0.108452809696
$ws.Range("R17").Value = 0.9760752872640001
$ws.Range("S17").Value = 0.0009767020798491111
$ws.Range("T17").Value = 0.0009767020798491111
